$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated get_alpha() results: the "_sell"/"t_*_sell" rows now reuse the prior
# "_buy"/"t_*_buy" row magnitudes (negated) once a negative monthly return is
# taken into account, and the "_buy", "t_*_buy" and "*_buy-sell" rows are
# recomputed (and filled in where previously blank) to match.

$ws.Range("B2").Value = 0.01947087286019047
$ws.Range("C2").Value = 0.01899497082536261
$ws.Range("D2").Value = 0.01824968330497252
$ws.Range("E2").Value = 0.01752827770277496

$ws.Range("B3").Value = 6.753759726330938
$ws.Range("C3").Value = 8.196272492707449
$ws.Range("D3").Value = 8.453782917330757
$ws.Range("E3").Value = 8.163224750770995

$ws.Range("B4").Value = -0.005514214498158991
$ws.Range("C4").Value = -0.007694010848738613
$ws.Range("D4").Value = -0.009578166933231158
$ws.Range("E4").Value = -0.01126444240632555

$ws.Range("B5").Value = -1.955155689973489
$ws.Range("C5").Value = -3.394831111512401
$ws.Range("D5").Value = -4.926115517268729
$ws.Range("E5").Value = -4.734337573985806

$ws.Range("B6").Value = 0.002861927112834738
$ws.Range("C6").Value = 0.003986331081739771
$ws.Range("D6").Value = 0.006700649650705811
$ws.Range("E6").Value = 0.005119225738973283

$ws.Range("B7").Value = 2.597917392994303
$ws.Range("C7").Value = 2.143298125438527
$ws.Range("D7").Value = 2.067360009852958
$ws.Range("E7").Value = 1.519453593188872

$ws.Range("B8").Value = 0.01764219997088727
$ws.Range("C8").Value = 0.01700852541705838
$ws.Range("D8").Value = 0.01764084704505658
$ws.Range("E8").Value = 0.01672176904611146

$ws.Range("B9").Value = 6.055472056248354
$ws.Range("C9").Value = 7.187689434648486
$ws.Range("D9").Value = 7.805666223097877
$ws.Range("E9").Value = 7.672869785858572

$ws.Range("B10").Value = -0.007657089661325975
$ws.Range("C10").Value = -0.008903211863784781
$ws.Range("D10").Value = -0.01063790674073911
$ws.Range("E10").Value = -0.01238057532362279

$ws.Range("B11").Value = -2.765221197521946
$ws.Range("C11").Value = -4.040882854193589
$ws.Range("D11").Value = -4.8014351120712
$ws.Range("E11").Value = -4.952453258095796

$ws.Range("B12").Value = 0.001904726305629638
$ws.Range("C12").Value = 0.001987492072762657
$ws.Range("D12").Value = 0.002724455621541621
$ws.Range("E12").Value = 0.004886627065307397

$ws.Range("B13").Value = 1.717983285850284
$ws.Range("C13").Value = 0.9590422325174918
$ws.Range("D13").Value = 1.052206627979532
$ws.Range("E13").Value = 1.265148738311833

$ws.Range("B14").Value = 0.01540097296011024
$ws.Range("C14").Value = 0.01567260764149626
$ws.Range("D14").Value = 0.01565745542650256
$ws.Range("E14").Value = 0.01500531424038237

$ws.Range("B15").Value = 5.181120564656332
$ws.Range("C15").Value = 6.3666298464954
$ws.Range("D15").Value = 6.6230269615689
$ws.Range("E15").Value = 6.531973301008986

$ws.Range("B16").Value = -0.007078685298556125
$ws.Range("C16").Value = -0.008680487604028096
$ws.Range("D16").Value = -0.01065621721080113
$ws.Range("E16").Value = -0.01259666415202041

$ws.Range("B17").Value = -2.515867263755929
$ws.Range("C17").Value = -3.804121870816258
$ws.Range("D17").Value = -4.619232663797836
$ws.Range("E17").Value = -5.124331086033146

$ws.Range("B18").Value = 0.001541359483814671
$ws.Range("C18").Value = 0.00215818822800695
$ws.Range("D18").Value = 0.005606844811557941
$ws.Range("E18").Value = 0.0008407733640093572

$ws.Range("B19").Value = 1.463944050248741
$ws.Range("C19").Value = 1.263018006722147
$ws.Range("D19").Value = 1.709327573824838
$ws.Range("E19").Value = 0.2129835129070369

$ws.Range("B20").Value = 0.01406294880386815
$ws.Range("C20").Value = 0.01411570727821992
$ws.Range("D20").Value = 0.0141174593798907
$ws.Range("E20").Value = 0.0139403844787541

$ws.Range("B21").Value = 4.56589033784384
$ws.Range("C21").Value = 5.357046784955265
$ws.Range("D21").Value = 5.807132996128346
$ws.Range("E21").Value = 5.804984948539256

$ws.Range("B22").Value = -0.008852954431593705
$ws.Range("C22").Value = -0.01150173881096942
$ws.Range("D22").Value = -0.01358360093365739
$ws.Range("E22").Value = -0.01447989462478875

$ws.Range("B23").Value = -3.123567132121082
$ws.Range("C23").Value = -4.75913123278851
$ws.Range("D23").Value = -5.571118917644933
$ws.Range("E23").Value = -5.632052754125775

$ws.Range("B24").Value = 0.0007886368958836297
$ws.Range("C24").Value = -0.0002990281703285538
$ws.Range("D24").Value = 0.0003324188366834002
$ws.Range("E24").Value = -0.003576121407434176

$ws.Range("B25").Value = 0.7412789507470753
$ws.Range("C25").Value = -0.1631020996807725
$ws.Range("D25").Value = 0.08618996605295239
$ws.Range("E25").Value = -0.8389275989676586
